$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2290.7778
$ws.Range("J28").Value = 6629.6
$ws.Range("L28").Value = 6629.6
$ws.Range("N28").Value = -7599.6

$ws.Range("H40").Value = 5403.8
$ws.Range("I40").Value = 5148.5713
$ws.Range("J40").Value = 5999.3335
$ws.Range("K40").Value = 5148.5713
$ws.Range("L40").Value = 5999.3335
$ws.Range("M40").Value = -4973.5713
$ws.Range("N40").Value = -6349.3335

$ws.Range("H62").Value = 2810.1667
$ws.Range("I62").Value = 1951
$ws.Range("J62").Value = 3669.3333
$ws.Range("K62").Value = 1951
$ws.Range("L62").Value = 3669.3333
$ws.Range("M62").Value = -1327
$ws.Range("N62").Value = -4917.3333

$ws.Range("H64").Value = 20839292
$ws.Range("I64").Value = 6040.61
$ws.Range("J64").Value = 142862620
$ws.Range("K64").Value = 6040.61
$ws.Range("L64").Value = 142862620
$ws.Range("M64").Value = -5792.61
$ws.Range("N64").Value = -142863116

$ws.Range("H65").Value = 2810.1667
$ws.Range("I65").Value = 1951
$ws.Range("J65").Value = 3669.3333
$ws.Range("K65").Value = 9755
$ws.Range("L65").Value = 18346.6665
$ws.Range("M65").Value = -6635
$ws.Range("N65").Value = -24586.6665

$ws.Range("H67").Value = 20839292
$ws.Range("I67").Value = 6040.61
$ws.Range("J67").Value = 142862620
$ws.Range("K67").Value = 6040.61
$ws.Range("L67").Value = 142862620
$ws.Range("M67").Value = -5182.61
$ws.Range("N67").Value = -142864336

$ws.Range("H112").Value = 2124.8
$ws.Range("J112").Value = 2124.8
$ws.Range("L112").Value = 6374.400000000001
$ws.Range("N112").Value = -8590.400000000001

$ws.Range("H118").Value = 723.2
$ws.Range("I118").Value = 703.7778
$ws.Range("J118").Value = 898
$ws.Range("K118").Value = 2111.3334
$ws.Range("L118").Value = 2694
$ws.Range("M118").Value = -454.3334
$ws.Range("N118").Value = -6008

$ws.Range("H125").Value = 1113
$ws.Range("I125").Value = 1203
$ws.Range("J125").Value = 1053
$ws.Range("K125").Value = 10827
$ws.Range("L125").Value = 9477
$ws.Range("M125").Value = -8367
$ws.Range("N125").Value = -14397

$ws.Range("H138").Value = 3375.1316
$ws.Range("J138").Value = 4192.077
$ws.Range("L138").Value = 12576.231
$ws.Range("N138").Value = -22856.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3322.125
$ws.Range("I32").Value = 3074
$ws.Range("K32").Value = 3074
$ws.Range("M32").Value = -2787

$ws.Range("H61").Value = 2426.756
$ws.Range("J61").Value = 3748.9167
$ws.Range("L61").Value = 3748.9167
$ws.Range("N61").Value = -4172.9167

$ws.Range("H63").Value = 2311.3845
$ws.Range("I63").Value = 2464.8
$ws.Range("K63").Value = 2464.8
$ws.Range("M63").Value = -1778.8

$ws.Range("H66").Value = 2311.3845
$ws.Range("I66").Value = 2464.8
$ws.Range("K66").Value = 12324
$ws.Range("M66").Value = -8892

$ws.Range("H88").Value = 3605.5454
$ws.Range("I88").Value = 1476.8
$ws.Range("J88").Value = 5379.5
$ws.Range("K88").Value = 1476.8
$ws.Range("L88").Value = 5379.5
$ws.Range("M88").Value = -1070.8
$ws.Range("N88").Value = -6191.5

$ws.Range("H91").Value = 3605.5454
$ws.Range("I91").Value = 1476.8
$ws.Range("J91").Value = 5379.5
$ws.Range("K91").Value = 1476.8
$ws.Range("L91").Value = 5379.5
$ws.Range("M91").Value = -72.79999999999995
$ws.Range("N91").Value = -8187.5

$ws.Range("H136").Value = 2426.756
$ws.Range("J136").Value = 3748.9167
$ws.Range("L136").Value = 11246.7501
$ws.Range("N136").Value = -16346.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 625.1111
$ws.Range("I22").Value = 516.5
$ws.Range("K22").Value = 516.5
$ws.Range("M22").Value = -343.5

$ws.Range("H134").Value = 3138.5
$ws.Range("I134").Value = 2911.353
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 8734.059000000001
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -6199.059000000001
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 792.3570999999999
$ws.Range("I22").Value = 1010.8889
$ws.Range("J22").Value = 399
$ws.Range("K22").Value = 1010.8889
$ws.Range("L22").Value = 399
$ws.Range("M22").Value = -660.8889
$ws.Range("N22").Value = -1099

$ws.Range("H31").Value = 3590.8723
$ws.Range("I31").Value = 2642.081
$ws.Range("K31").Value = 2642.081
$ws.Range("M31").Value = -2347.081

$ws.Range("H34").Value = 3590.8723
$ws.Range("I34").Value = 2642.081
$ws.Range("K34").Value = 2642.081
$ws.Range("M34").Value = -2440.081

$ws.Range("H58").Value = 2738.2727
$ws.Range("I58").Value = 1424.2
$ws.Range("J58").Value = 3833.3333
$ws.Range("K58").Value = 1424.2
$ws.Range("L58").Value = 3833.3333
$ws.Range("M58").Value = -1221.2
$ws.Range("N58").Value = -4239.3333

$ws.Range("H107").Value = 3571694.2
$ws.Range("J107").Value = 788
$ws.Range("L107").Value = 788
$ws.Range("N107").Value = -4628

$ws.Range("H131").Value = 58442.332
$ws.Range("J131").Value = 58442.332
$ws.Range("L131").Value = 58442.332
$ws.Range("N131").Value = -68522.33199999999

$ws.Range("H134").Value = 3034.8096
$ws.Range("I134").Value = 2837.1177
$ws.Range("K134").Value = 8511.3531
$ws.Range("M134").Value = -5976.3531

$ws.Range("H136").Value = 2738.2727
$ws.Range("I136").Value = 1424.2
$ws.Range("J136").Value = 3833.3333
$ws.Range("K136").Value = 4272.6
$ws.Range("L136").Value = 11499.9999
$ws.Range("M136").Value = -1722.6
$ws.Range("N136").Value = -16599.9999

$ws.Range("H139").Value = 55221
$ws.Range("J139").Value = 55221
$ws.Range("L139").Value = 55221
$ws.Range("N139").Value = -65501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 311.08334
$ws.Range("J107").Value = 311.08334
$ws.Range("L107").Value = 933.2500200000001
$ws.Range("N107").Value = -4773.25002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 160291
$ws.Range("I70").Value = 290255.28
$ws.Range("J70").Value = 8666
$ws.Range("K70").Value = 290255.28
$ws.Range("L70").Value = 8666
$ws.Range("M70").Value = -289985.28
$ws.Range("N70").Value = -9206

$ws.Range("H73").Value = 160291
$ws.Range("I73").Value = 290255.28
$ws.Range("J73").Value = 8666
$ws.Range("K73").Value = 290255.28
$ws.Range("L73").Value = 8666
$ws.Range("M73").Value = -289319.28
$ws.Range("N73").Value = -10538

$ws.Range("H122").Value = 4482.1333
$ws.Range("I122").Value = 3633.2307
$ws.Range("K122").Value = 10899.6921
$ws.Range("M122").Value = -8449.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 17985
$ws.Range("I58").Value = 17985
$ws.Range("K58").Value = 17985
$ws.Range("M58").Value = -17677

$ws.Range("H112").Value = 48619
$ws.Range("J112").Value = 48619
$ws.Range("L112").Value = 48619
$ws.Range("N112").Value = -51573

$ws.Range("H123").Value = 92248.25
$ws.Range("J123").Value = 92248.25
$ws.Range("L123").Value = 92248.25
$ws.Range("N123").Value = -102048.25

$ws.Range("H126").Value = 3306.9285
$ws.Range("I126").Value = 3307.4614
$ws.Range("K126").Value = 9922.3842
$ws.Range("M126").Value = -7452.3842

$ws.Range("H136").Value = 47502.5
$ws.Range("J136").Value = 47502.5
$ws.Range("L136").Value = 142507.5
$ws.Range("N136").Value = -147607.5

Write-Output "Applied Gilgamesh_Profits market-data refresh across 8 sheets."
